$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.88649582862854
$ws.Range("B1").Value = 1.243007779121399
$ws.Range("C1").Value = 1.767964959144592
$ws.Range("D1").Value = 5.205566883087158
$ws.Range("E1").Value = 2.073575258255005
